$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DynamicListingPage")

# New "ExpectedFilterInfo" column (C) with header + values for each existing
# test row, describing the filter breadcrumb expected for that Path/ContentType.
$ws.Range("C1").Value = "ExpectedFilterInfo"
$ws.Range("C2").Value = "disease|breast-cancer|none|none|"
$ws.Range("C3").Value = "disease|breast-cancer|treatment|none|"
$ws.Range("C4").Value = "disease|breast-cancer|treatment|trastuzumab|"
$ws.Range("C5").Value = "intervention|trastuzumab|none|"
$ws.Range("C6").Value = "intervention|trastuzumab|treatment|"
$ws.Range("C7").Value = "manual parameters|"

# Match the bold/filled header formatting already used for A1:B1.
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Size the new column similarly to the other bestFit columns.
$ws.Columns.Item(3).ColumnWidth = 44
